# Updating the Forecast Portfolio
# Shift every timestamp in column A (rows 2-97) forward by 16 days, and
# replace the "Actual Production (MW)" values in column B with the new
# forecast readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Actual Production (MW) values for rows 2..97 (96 values).
$newB = @(
    908,949,1066,1179,1234,1258,1298,1310,1376,1354,
    1296,1254,1239,1262,1214,1137,1161,1254,1365,1443,
    1442,1468,1483,0,1583,1700,1790,1876,1890,1813,
    1830,1819,1762,1634,1503,1441,1428,1467,1504,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0
)

for ($i = 0; $i -lt 96; $i++) {
    $row = $i + 2

    # Shift the timestamp forward by 16 days (keeps the same time-of-day).
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value2 = $dateCell.Value2 + 16

    # Write the new production figure.
    $ws.Cells.Item($row, 2).Value2 = $newB[$i]
}
